$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 corresponds to the Schottky diode BOM line (D, 40V 3A Schottky, B340A-13-F DII).
# The power supply input now uses a second diode (D44) in addition to D1,
# so the quantity (Anzahl) goes from 1 to 2 and the positions (Positionen)
# list is updated from "D1" to "D1, D44".
$ws.Range("A6").Value = 2
$ws.Range("D6").Value = "D1, D44"

# Keep the active selection consistent with the edited area, as seen in the diff.
$ws.Range("D10").Select()
